# Update the "F" column (观展人数/浏览量 style numeric counter) figures on the
# "展览" sheet and the "全部类型" sheet to reflect newly generated output
# (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) -------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value  = 3247
$ws1.Range("F5").Value  = 2321
$ws1.Range("F8").Value  = 1317
$ws1.Range("F9").Value  = 1057
$ws1.Range("F10").Value = 270
$ws1.Range("F11").Value = 488
$ws1.Range("F12").Value = 1166
$ws1.Range("F14").Value = 88
$ws1.Range("F15").Value = 539
$ws1.Range("F16").Value = 8160
$ws1.Range("F17").Value = 356
$ws1.Range("F19").Value = 230
$ws1.Range("F20").Value = 243
$ws1.Range("F23").Value = 555
$ws1.Range("F25").Value = 1145
$ws1.Range("F27").Value = 1889
$ws1.Range("F28").Value = 928
$ws1.Range("F30").Value = 1707
$ws1.Range("F32").Value = 1910
$ws1.Range("F34").Value = 13
$ws1.Range("F35").Value = 6
$ws1.Range("F36").Value = 60
$ws1.Range("F39").Value = 48
$ws1.Range("F40").Value = 202
$ws1.Range("F41").Value = 365
$ws1.Range("F43").Value = 238

# --- Sheet "全部类型" (sheet4) ---------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value  = 3247
$ws4.Range("F7").Value  = 2321
$ws4.Range("F10").Value = 1317
$ws4.Range("F12").Value = 1057
$ws4.Range("F13").Value = 270
$ws4.Range("F14").Value = 488
$ws4.Range("F15").Value = 88
$ws4.Range("F16").Value = 539
$ws4.Range("F17").Value = 8160
$ws4.Range("F18").Value = 356
$ws4.Range("F21").Value = 230
$ws4.Range("F22").Value = 243
$ws4.Range("F25").Value = 555
$ws4.Range("F27").Value = 1145
$ws4.Range("F29").Value = 1889
$ws4.Range("F30").Value = 928
$ws4.Range("F32").Value = 1707
$ws4.Range("F33").Value = 1910
$ws4.Range("F35").Value = 13
$ws4.Range("F36").Value = 6
$ws4.Range("F37").Value = 60
$ws4.Range("F40").Value = 48
$ws4.Range("F41").Value = 202
$ws4.Range("F42").Value = 365
$ws4.Range("F49").Value = 238
